# Regenerate orders with updated distance/size codes.
# Mapping: D51 -> D55, D64 -> D69, D80 -> D86, S30 -> S31
# Applied as a global find/replace over every cell on the active sheet,
# which naturally covers the Condition, Filename_Left, Filename_Right,
# Distance and Size columns (all of which embed these tokens).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("D51", "D55")
$ws.Cells.Replace("D64", "D69")
$ws.Cells.Replace("D80", "D86")
$ws.Cells.Replace("S30", "S31")
